# Realestate Update resale numbers 2023-06-08 09:37
# Append a new data row (row 29) to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 29

# Columns A-D hold text in this sheet (dates/times/weekday/week stored as
# literal strings, not real dates/numbers) - force text format so Excel
# doesn't auto-convert "2023-06-08" / "23" into a date serial / number,
# then clear the formatting again so the new row keeps the sheet's plain
# (unstyled) look, matching the other data rows.
$ws.Range("A$row`:D$row").NumberFormat = "@"

$ws.Cells.Item($row, 1).Value  = "2023-06-08"
$ws.Cells.Item($row, 2).Value  = "09:36:06"
$ws.Cells.Item($row, 3).Value  = "Thursday"
$ws.Cells.Item($row, 4).Value  = "23"

$ws.Range("A$row`:D$row").ClearFormats()

$ws.Cells.Item($row, 5).Value  = 117405
$ws.Cells.Item($row, 6).Value  = 134379
$ws.Cells.Item($row, 7).Value  = 159895
$ws.Cells.Item($row, 8).Value  = 130756
$ws.Cells.Item($row, 9).Value  = 175449
$ws.Cells.Item($row, 10).Value = 112786
$ws.Cells.Item($row, 11).Value = 200877
$ws.Cells.Item($row, 12).Value = 220890
$ws.Cells.Item($row, 13).Value = 172742
$ws.Cells.Item($row, 14).Value = 119955
$ws.Cells.Item($row, 15).Value = 38566
$ws.Cells.Item($row, 16).Value = 34485
$ws.Cells.Item($row, 17).Value = 50734
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36422
$ws.Cells.Item($row, 20).Value = -1
